$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching style of existing headers
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add data values for the new column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
